$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 3)
$ws.Range("A3").Value = "NM7"
$ws.Range("D3").Value = "New Subproduct2"

# C3 should carry the same "class" style/border as C2 (the bordered "Usage" style),
# so copy C2's formatting down into C3 before setting its value.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C3").Value = "SAP & E1"

$ws.Range("D4").Select()
